# Automatic update of files.
# Bump the "Förändrad" (Changed) date in column C by one day for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlUp = -4162

# Column C = "Förändrad" (header is in row 1, data starts at row 2)
$colChanged = 3
$headerRow = 1

$lastRow = $ws.Cells.Item($ws.Rows.Count, $colChanged).End($xlUp).Row
if ($lastRow -lt $headerRow) { $lastRow = $headerRow }

for ($r = $headerRow + 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $colChanged)
    $current = $cell.Value2
    if ($current -ne $null) {
        $cell.Value2 = $current + 1
    }
}
